$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.5403
$ws.Range("D3").Value = -6.79189999999999
$ws.Range("E8").Value = 16.59020000000001
$ws.Range("E11").Value = 16.8537
$ws.Range("B12").Value = 4.952299999999997
$ws.Range("C14").Value = -13.7595
$ws.Range("E14").Value = 16.88640000000001
$ws.Range("E15").Value = 16.1896
$ws.Range("E17").Value = 16.82030000000002
$ws.Range("D20").Value = -7.502999999999998
$ws.Range("D25").Value = -7.601100000000006
$ws.Range("C26").Value = -12.1164
$ws.Range("E26").Value = 15.4068
$ws.Range("B27").Value = 6.040800000000004
$ws.Range("D30").Value = -7.514600000000004
$ws.Range("C31").Value = -13.6821
$ws.Range("B32").Value = 6.503999999999998
$ws.Range("C35").Value = -12.36740000000001
$ws.Range("B36").Value = 9.094299999999999
$ws.Range("E36").Value = 16.4708
$ws.Range("C37").Value = -13.85259999999999
$ws.Range("B38").Value = 5.687299999999999
$ws.Range("D44").Value = -7.209200000000001
$ws.Range("C45").Value = -13.30989999999999
$ws.Range("B46").Value = 6.8312
$ws.Range("D47").Value = -7.290799999999999
$ws.Range("C52").Value = -10.9427
$ws.Range("B54").Value = 4.926000000000001
$ws.Range("B55").Value = 5.270499999999998
$ws.Range("B56").Value = 4.5176
$ws.Range("C57").Value = -14.62829999999999
$ws.Range("D58").Value = -7.760800000000002
$ws.Range("E64").Value = 17.5124
$ws.Range("B67").Value = 4.891499999999994
$ws.Range("B69").Value = 5.381599999999993
$ws.Range("B72").Value = 5.980899999999997
$ws.Range("D78").Value = -7.901000000000003
$ws.Range("E79").Value = 17.94270000000002
$ws.Range("C81").Value = -13.2573
$ws.Range("B83").Value = 5.446099999999996
$ws.Range("C83").Value = -13.6373
$ws.Range("D84").Value = -8.412799999999999
$ws.Range("B86").Value = 4.752000000000002
$ws.Range("D89").Value = -6.469499999999995
$ws.Range("E89").Value = 18.06830000000001
$ws.Range("B91").Value = 5.189499999999999
$ws.Range("D91").Value = -6.524099999999997
$ws.Range("D92").Value = -6.361400000000004
$ws.Range("B93").Value = 6.565
$ws.Range("D96").Value = -7.555500000000003
$ws.Range("B99").Value = 4.589099999999999
$ws.Range("C100").Value = -12.8598
$ws.Range("C102").Value = -13.8969
$ws.Range("D102").Value = -7.831499999999997
